# Updates cryptos list cell values (Price and Volume(1h)) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.375.46"
$ws.Range("E2").Value = "  -2.10%  "
$ws.Range("D3").Value = "2.624.82"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.92"
$ws.Range("E5").Value = "  -3.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.39"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.649"
$ws.Range("E7").Value = "  +5.74%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -5.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.77"
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.388"
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.39"
$ws.Range("E13").Value = "  -2.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000185"
$ws.Range("E14").Value = "  -6.90%  "
$ws.Range("D15").Value = "3.097.12"
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("D16").Value = "64.262.64"
$ws.Range("E16").Value = "  -2.07%  "
$ws.Range("D17").Value = "2.625.98"
$ws.Range("E17").Value = "  -2.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.20"
$ws.Range("E18").Value = "  -4.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.66"
$ws.Range("E19").Value = "  -2.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.34"
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "344.36"
$ws.Range("E21").Value = "  -2.06%  "
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.61"
$ws.Range("E23").Value = "  -2.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000112"
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("E25").Value = "  +2.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.29"
$ws.Range("E26").Value = "  -4.09%  "
$ws.Range("E27").Value = "  -3.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "553.21"
$ws.Range("E28").Value = "  +3.83%  "
$ws.Range("E29").Value = "  -2.50%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.88"
$ws.Range("E31").Value = "  -1.58%  "
$ws.Range("E32").Value = "  -3.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.72"
$ws.Range("E33").Value = "  -3.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.38"
$ws.Range("E34").Value = "  -1.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.24"
$ws.Range("E35").Value = "  -4.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.411"
$ws.Range("E36").Value = "  -2.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.91"
$ws.Range("E37").Value = "  -3.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.91"
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "153.37"
$ws.Range("E40").Value = "  -2.48%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.42"
$ws.Range("E42").Value = "  +3.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "157.82"
$ws.Range("E43").Value = "  -3.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.96"
$ws.Range("E44").Value = "  -3.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0597"
$ws.Range("E45").Value = "  -2.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.61"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.632"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("E48").Value = "  +2.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0249"
$ws.Range("E49").Value = "  -3.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.98"
$ws.Range("E50").Value = "  -5.08%  "
$ws.Range("E51").Value = "  -5.98%  "
